# Applies numeric corrections to the Kujata_Profits price/profit columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets, per the scheduled-runner pricing refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1833.7966
$ws.Range("I98").Value = 1927.9259
$ws.Range("K98").Value = 1927.9259
$ws.Range("M98").Value = -429.9259
$ws.Range("H112").Value = 2543.1304
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 2608.7273
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 7826.1819
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -10042.1819
$ws.Range("H115").Value = 435
$ws.Range("I115").Value = 435
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1305
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 262
$ws.Range("N115").ClearContents()
$ws.Range("H122").Value = 1833.7966
$ws.Range("I122").Value = 1927.9259
$ws.Range("K122").Value = 5783.7777
$ws.Range("M122").Value = -3333.7777
$ws.Range("H132").Value = 9264175
$ws.Range("I132").Value = 10421939
$ws.Range("J132").Value = 2056.5
$ws.Range("K132").Value = 31265817
$ws.Range("L132").Value = 6169.5
$ws.Range("M132").Value = -31263287
$ws.Range("N132").Value = -11229.5
$ws.Range("H133").Value = 32913.625
$ws.Range("J133").Value = 32913.625
$ws.Range("L133").Value = 32913.625
$ws.Range("N133").Value = -43033.625
$ws.Range("H137").Value = 1500.9736
$ws.Range("I137").Value = 996.3158
$ws.Range("J137").Value = 2005.6316
$ws.Range("K137").Value = 2988.9474
$ws.Range("L137").Value = 6016.8948
$ws.Range("M137").Value = -438.9474
$ws.Range("N137").Value = -11116.8948
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4615
$ws.Range("I31").Value = 4615
$ws.Range("K31").Value = 4615
$ws.Range("M31").Value = -4321
$ws.Range("H32").Value = 9985.841
$ws.Range("I32").Value = 7694.6855
$ws.Range("J32").Value = 18895.889
$ws.Range("K32").Value = 7694.6855
$ws.Range("L32").Value = 18895.889
$ws.Range("M32").Value = -7407.6855
$ws.Range("N32").Value = -19469.889
$ws.Range("H74").Value = 1182.3077
$ws.Range("I74").Value = 979.0909
$ws.Range("J74").Value = 2300
$ws.Range("K74").Value = 979.0909
$ws.Range("L74").Value = 2300
$ws.Range("M74").Value = -105.0909
$ws.Range("N74").Value = -4048
$ws.Range("H77").Value = 1182.3077
$ws.Range("I77").Value = 979.0909
$ws.Range("J77").Value = 2300
$ws.Range("K77").Value = 4895.4545
$ws.Range("L77").Value = 11500
$ws.Range("M77").Value = -527.4544999999998
$ws.Range("N77").Value = -20236
$ws.Range("H102").Value = 15153659
$ws.Range("I102").Value = 18520476
$ws.Range("K102").Value = 18520476
$ws.Range("M102").Value = -18518854
$ws.Range("H132").Value = 2694.6316
$ws.Range("I132").Value = 2428.6155
$ws.Range("K132").Value = 7285.8465
$ws.Range("M132").Value = -4755.8465
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 14780
$ws.Range("J81").Value = 14780
$ws.Range("L81").Value = 14780
$ws.Range("N81").Value = -16902
$ws.Range("H84").Value = 14780
$ws.Range("J84").Value = 14780
$ws.Range("L84").Value = 44340
$ws.Range("N84").Value = -54948
$ws.Range("H134").Value = 2028.4231
$ws.Range("I134").Value = 1518.8636
$ws.Range("K134").Value = 4556.5908
$ws.Range("M134").Value = -2021.5908
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1575.3846
$ws.Range("I31").Value = 1189.0952
$ws.Range("J31").Value = 3197.8
$ws.Range("K31").Value = 1189.0952
$ws.Range("L31").Value = 3197.8
$ws.Range("M31").Value = -894.0952
$ws.Range("N31").Value = -3787.8
$ws.Range("H34").Value = 1575.3846
$ws.Range("I34").Value = 1189.0952
$ws.Range("J34").Value = 3197.8
$ws.Range("K34").Value = 1189.0952
$ws.Range("L34").Value = 3197.8
$ws.Range("M34").Value = -987.0952
$ws.Range("N34").Value = -3601.8
$ws.Range("H94").Value = 1285.25
$ws.Range("I94").Value = 1175.4
$ws.Range("J94").Value = 1363.7142
$ws.Range("K94").Value = 1175.4
$ws.Range("L94").Value = 1363.7142
$ws.Range("M94").Value = -724.4000000000001
$ws.Range("N94").Value = -2265.7142
$ws.Range("H132").Value = 6196.625
$ws.Range("I132").Value = 7558.1875
$ws.Range("J132").Value = 3473.5
$ws.Range("K132").Value = 22674.5625
$ws.Range("L132").Value = 10420.5
$ws.Range("M132").Value = -20144.5625
$ws.Range("N132").Value = -15480.5
$ws.Range("H134").Value = 16130802
$ws.Range("I134").Value = 1798.36
$ws.Range("K134").Value = 5395.08
$ws.Range("M134").Value = -2860.08
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 645.381
$ws.Range("I113").Value = 499
$ws.Range("J113").Value = 679.82355
$ws.Range("K113").Value = 1497
$ws.Range("L113").Value = 2039.47065
$ws.Range("M113").Value = 673
$ws.Range("N113").Value = -6379.470649999999
$ws.Range("H131").Value = 27782128
$ws.Range("J131").Value = 5163.3335
$ws.Range("L131").Value = 15490.0005
$ws.Range("N131").Value = -25570.0005
$ws.Range("H132").Value = 1103.0769
$ws.Range("I132").Value = 1095.5555
$ws.Range("J132").Value = 1120
$ws.Range("K132").Value = 9859.9995
$ws.Range("L132").Value = 10080
$ws.Range("M132").Value = -7329.9995
$ws.Range("N132").Value = -15140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4625.1
$ws.Range("I46").Value = 617
$ws.Range("J46").Value = 6342.857
$ws.Range("K46").Value = 617
$ws.Range("L46").Value = 6342.857
$ws.Range("M46").Value = -429
$ws.Range("N46").Value = -6718.857
$ws.Range("H132").Value = 103255.09
$ws.Range("I132").Value = 37101
$ws.Range("J132").Value = 128062.875
$ws.Range("K132").Value = 111303
$ws.Range("L132").Value = 384188.625
$ws.Range("M132").Value = -108773
$ws.Range("N132").Value = -389248.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 10041
$ws.Range("J44").Value = 10041
$ws.Range("L44").Value = 10041
$ws.Range("N44").Value = -11149
$ws.Range("H63").Value = 12709.857
$ws.Range("J63").Value = 15548.6
$ws.Range("L63").Value = 15548.6
$ws.Range("N63").Value = -16796.6
$ws.Range("H66").Value = 12709.857
$ws.Range("J66").Value = 15548.6
$ws.Range("L66").Value = 46645.8
$ws.Range("N66").Value = -52885.8
$ws.Range("H96").Value = 2799.7856
$ws.Range("I96").Value = 1809.3
$ws.Range("K96").Value = 1809.3
$ws.Range("M96").Value = -436.3
$ws.Range("H132").Value = 4076.2942
$ws.Range("I132").Value = 5259.8
$ws.Range("K132").Value = 15779.4
$ws.Range("M132").Value = -13249.4
$ws.Range("H136").Value = 961.7273
$ws.Range("I136").Value = 936.6818
$ws.Range("J136").Value = 1011.8182
$ws.Range("K136").Value = 2810.0454
$ws.Range("L136").Value = 3035.4546
$ws.Range("M136").Value = -260.0454
$ws.Range("N136").Value = -8135.4546
